$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: "ook after:" -> "look after:" (new leading run containing "l",
# keeping the original "ook after:" run separate).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lookAfterPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "ook after:`r") {
        $lookAfterPara = $p
    }
}

if ($lookAfterPara -ne $null) {
    $xml1 = "<w:p $wns><w:r><w:t>l</w:t></w:r><w:r><w:t>ook after:</w:t></w:r></w:p>"
    $lookAfterPara.Range.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Change 2: after the "minimum nights:" paragraph, insert the new analysis /
# RMSE notes block.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$minNightsPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "minimum nights:`r") {
        $minNightsPara = $p
    }
}

if ($minNightsPara -ne $null) {
    $parts = @()
    $parts += "<w:p $wns><w:r><w:t>minimum nights:</w:t></w:r></w:p>"
    $parts += "<w:p $wns/>"
    $parts += "<w:p $wns><w:r><w:t>For analyis:</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>t1 -&gt; model table for cross validation</w:t></w:r><w:r><w:t xml:space=`"preserve`"> and models</w:t></w:r></w:p>"
    $parts += "<w:p $wns/>"
    $parts += "<w:p $wns/>"
    $parts += "<w:p $wns><w:r><w:t>RMSE:</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>&gt; model3_level_holdout_rmse</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>[1] 30</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>&gt; model3_level_work_rmse</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>[1] 28.3</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>&gt;</w:t></w:r></w:p>"
    $parts += "<w:p $wns/>"
    $parts += "<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t>&gt; model4_level_holdout_rmse</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>[1] 30</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>&gt; model4_level_work_rmse</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>[1] 28.2</w:t></w:r></w:p>"
    $parts += "<w:p $wns><w:r><w:t>&gt;</w:t></w:r></w:p>"
    $parts += "<w:p $wns/>"

    $xml2 = [string]::Join("", $parts)
    $minNightsPara.Range.InsertXML($xml2)
}

Write-Output "edit complete"
